# Auto-generated edit script applying scheduled market-data refresh to Chocobo_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 25642728
$ws.Range("J112").Value = 1781.6389
$ws.Range("L112").Value = 5344.9167
$ws.Range("N112").Value = -7560.9167
$ws.Range("H113").Value = 7294.737
$ws.Range("I113").Value = 4597
$ws.Range("J113").Value = 8014.1333
$ws.Range("K113").Value = 4597
$ws.Range("L113").Value = 8014.1333
$ws.Range("M113").Value = -1343
$ws.Range("N113").Value = -14522.1333
$ws.Range("H129").Value = 838.84375
$ws.Range("J129").Value = 961.451
$ws.Range("L129").Value = 2884.353
$ws.Range("N129").Value = -12884.353
$ws.Range("H132").Value = 33340874
$ws.Range("I132").Value = 47626972
$ws.Range("J132").Value = 6645.5557
$ws.Range("K132").Value = 142880916
$ws.Range("L132").Value = 19936.6671
$ws.Range("M132").Value = -142878386
$ws.Range("N132").Value = -24996.6671
$ws.Range("H137").Value = 2705.4883
$ws.Range("I137").Value = 1195.4138
$ws.Range("K137").Value = 3586.2414
$ws.Range("M137").Value = -1036.2414
$ws.Range("H138").Value = 4728.64
$ws.Range("I138").Value = 685.35297
$ws.Range("J138").Value = 5556.783
$ws.Range("K138").Value = 2056.05891
$ws.Range("L138").Value = 16670.349
$ws.Range("M138").Value = 3083.94109
$ws.Range("N138").Value = -26950.349
$ws.Range("H141").Value = 3717.4
$ws.Range("I141").Value = 3657.6052
$ws.Range("J141").Value = 4042
$ws.Range("K141").Value = 10972.8156
$ws.Range("L141").Value = 12126
$ws.Range("M141").Value = -5792.8156
$ws.Range("N141").Value = -22486

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3331.647
$ws.Range("I74").Value = 3409.1853
$ws.Range("K74").Value = 3409.1853
$ws.Range("M74").Value = -2535.1853
$ws.Range("H76").Value = 26825.092
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 26825.092
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 26825.092
$ws.Range("M76").Value = $null
$ws.Range("N76").Value = -27501.092
$ws.Range("H77").Value = 3331.647
$ws.Range("I77").Value = 3409.1853
$ws.Range("K77").Value = 17045.9265
$ws.Range("M77").Value = -12677.9265
$ws.Range("H79").Value = 26825.092
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 26825.092
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 26825.092
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -29165.092
$ws.Range("H115").Value = 29890
$ws.Range("J115").Value = 29890
$ws.Range("L115").Value = 29890
$ws.Range("N115").Value = -33024
$ws.Range("H132").Value = 2113.7886
$ws.Range("I132").Value = 1300.5128
$ws.Range("J132").Value = 4553.615
$ws.Range("K132").Value = 3901.5384
$ws.Range("L132").Value = 13660.845
$ws.Range("M132").Value = -1371.5384
$ws.Range("N132").Value = -18720.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1683.4023
$ws.Range("I134").Value = 1057.5555
$ws.Range("K134").Value = 3172.6665
$ws.Range("M134").Value = -637.6664999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7814971.5
$ws.Range("I31").Value = 1361.5111
$ws.Range("J31").Value = 26320890
$ws.Range("K31").Value = 1361.5111
$ws.Range("L31").Value = 26320890
$ws.Range("M31").Value = -1066.5111
$ws.Range("N31").Value = -26321480
$ws.Range("H34").Value = 7814971.5
$ws.Range("I34").Value = 1361.5111
$ws.Range("J34").Value = 26320890
$ws.Range("K34").Value = 1361.5111
$ws.Range("L34").Value = 26320890
$ws.Range("M34").Value = -1159.5111
$ws.Range("N34").Value = -26321294
$ws.Range("H58").Value = 1206.49
$ws.Range("I58").Value = 1330.9487
$ws.Range("J58").Value = 765.2273
$ws.Range("K58").Value = 1330.9487
$ws.Range("L58").Value = 765.2273
$ws.Range("M58").Value = -1127.9487
$ws.Range("N58").Value = -1171.2273
$ws.Range("H99").Value = 6671055.5
$ws.Range("I99").Value = 10529631
$ws.Range("J99").Value = 6242.636
$ws.Range("K99").Value = 10529631
$ws.Range("L99").Value = 6242.636
$ws.Range("M99").Value = -10528133
$ws.Range("N99").Value = -9238.636
$ws.Range("H105").Value = 2299.818
$ws.Range("I105").Value = 2102.5
$ws.Range("J105").Value = 2826
$ws.Range("K105").Value = 2102.5
$ws.Range("L105").Value = 2826
$ws.Range("M105").Value = -355.5
$ws.Range("N105").Value = -6320
$ws.Range("H126").Value = 6671055.5
$ws.Range("I126").Value = 10529631
$ws.Range("J126").Value = 6242.636
$ws.Range("K126").Value = 31588893
$ws.Range("L126").Value = 18727.908
$ws.Range("M126").Value = -31586423
$ws.Range("N126").Value = -23667.908
$ws.Range("H132").Value = 2022.9215
$ws.Range("I132").Value = 1646.7561
$ws.Range("J132").Value = 3565.2
$ws.Range("K132").Value = 4940.2683
$ws.Range("L132").Value = 10695.6
$ws.Range("M132").Value = -2410.2683
$ws.Range("N132").Value = -15755.6
$ws.Range("H134").Value = 3047.492
$ws.Range("I134").Value = 3587.2778
$ws.Range("J134").Value = 2327.7778
$ws.Range("K134").Value = 10761.8334
$ws.Range("L134").Value = 6983.3334
$ws.Range("M134").Value = -8226.8334
$ws.Range("N134").Value = -12053.3334
$ws.Range("H136").Value = 1206.49
$ws.Range("I136").Value = 1330.9487
$ws.Range("J136").Value = 765.2273
$ws.Range("K136").Value = 3992.8461
$ws.Range("L136").Value = 2295.6819
$ws.Range("M136").Value = -1442.8461
$ws.Range("N136").Value = -7395.6819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 606.92725
$ws.Range("I113").Value = 505.3095
$ws.Range("K113").Value = 1515.9285
$ws.Range("M113").Value = 654.0715
$ws.Range("H131").Value = 844.6061
$ws.Range("J131").Value = 992.3913
$ws.Range("L131").Value = 2977.1739
$ws.Range("N131").Value = -13057.1739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6652.758
$ws.Range("I70").Value = 6051.7144
$ws.Range("K70").Value = 6051.7144
$ws.Range("M70").Value = -5781.7144
$ws.Range("H73").Value = 6652.758
$ws.Range("I73").Value = 6051.7144
$ws.Range("K73").Value = 6051.7144
$ws.Range("M73").Value = -5115.7144
$ws.Range("H102").Value = 1869.9736
$ws.Range("I102").Value = 1614.8948
$ws.Range("J102").Value = 2125.0527
$ws.Range("K102").Value = 1614.8948
$ws.Range("L102").Value = 2125.0527
$ws.Range("M102").Value = 7.105199999999968
$ws.Range("N102").Value = -5369.0527
$ws.Range("H132").Value = 2358.98
$ws.Range("I132").Value = 1520.5938
$ws.Range("J132").Value = 3849.4443
$ws.Range("K132").Value = 4561.7814
$ws.Range("L132").Value = 11548.3329
$ws.Range("M132").Value = -2031.7814
$ws.Range("N132").Value = -16608.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4800.5
$ws.Range("I7").Value = 3488.6667
$ws.Range("J7").Value = 6487.143
$ws.Range("K7").Value = 3488.6667
$ws.Range("L7").Value = 6487.143
$ws.Range("M7").Value = -3376.6667
$ws.Range("N7").Value = -6711.143
$ws.Range("H40").Value = 5373.091
$ws.Range("I40").Value = 3755.9333
$ws.Range("J40").Value = 8838.429
$ws.Range("K40").Value = 3755.9333
$ws.Range("L40").Value = 8838.429
$ws.Range("M40").Value = -3619.9333
$ws.Range("N40").Value = -9110.429
$ws.Range("H74").Value = 34776.668
$ws.Range("I74").Value = 17000
$ws.Range("J74").Value = 39855.715
$ws.Range("K74").Value = 17000
$ws.Range("L74").Value = 39855.715
$ws.Range("M74").Value = -16002
$ws.Range("N74").Value = -41851.715
$ws.Range("H77").Value = 34776.668
$ws.Range("I77").Value = 17000
$ws.Range("J77").Value = 39855.715
$ws.Range("K77").Value = 51000
$ws.Range("L77").Value = 119567.145
$ws.Range("M77").Value = -46008
$ws.Range("N77").Value = -129551.145
$ws.Range("H122").Value = 6686.0713
$ws.Range("I122").Value = 3650
$ws.Range("J122").Value = 7900.5
$ws.Range("K122").Value = 10950
$ws.Range("L122").Value = 23701.5
$ws.Range("M122").Value = -8500
$ws.Range("N122").Value = -28601.5
$ws.Range("H126").Value = 4800.5
$ws.Range("I126").Value = 3488.6667
$ws.Range("J126").Value = 6487.143
$ws.Range("K126").Value = 10466.0001
$ws.Range("L126").Value = 19461.429
$ws.Range("M126").Value = -7996.000100000001
$ws.Range("N126").Value = -24401.429
$ws.Range("H132").Value = 3384.125
$ws.Range("I132").Value = 1195.2325
$ws.Range("J132").Value = 6629.724
$ws.Range("K132").Value = 3585.6975
$ws.Range("L132").Value = 19889.172
$ws.Range("M132").Value = -1055.6975
$ws.Range("N132").Value = -24949.172
$ws.Range("H136").Value = 2551.7234
$ws.Range("I136").Value = 1677.2188
$ws.Range("J136").Value = 4417.3335
$ws.Range("K136").Value = 5031.6564
$ws.Range("L136").Value = 13252.0005
$ws.Range("M136").Value = -2481.6564
$ws.Range("N136").Value = -18352.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2905.2
$ws.Range("I122").Value = 1784.1765
$ws.Range("J122").Value = 5287.375
$ws.Range("K122").Value = 5352.529500000001
$ws.Range("L122").Value = 15862.125
$ws.Range("M122").Value = -2902.529500000001
$ws.Range("N122").Value = -20762.125
$ws.Range("H126").Value = 238426.98
$ws.Range("I126").Value = 1398.2916
$ws.Range("K126").Value = 4194.8748
$ws.Range("M126").Value = -1724.8748
$ws.Range("H132").Value = 4505675
$ws.Range("I132").Value = 548.3148
$ws.Range("J132").Value = 16669517
$ws.Range("K132").Value = 1644.9444
$ws.Range("L132").Value = 50008551
$ws.Range("M132").Value = 885.0556000000001
$ws.Range("N132").Value = -50013611
$ws.Range("H136").Value = 2012.0333
$ws.Range("I136").Value = 728.2222
$ws.Range("J136").Value = 3937.75
$ws.Range("K136").Value = 2012.0333
$ws.Range("L136").Value = 11813.25
$ws.Range("M136").Value = 365.3334
$ws.Range("N136").Value = -16913.25
